# Swap the data values (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, and the weighted-price helper column) between
# rows 4-5 and rows 6-7, leaving all other columns (A, B, C, E-L, Q, R, T)
# untouched, as described by the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <- old row 6 values
$ws.Range("D4").Value = 44223
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 3500
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3750
$ws.Range("S4").Value = 1875

# Row 5 <- old row 7 values
$ws.Range("D5").Value = 44223
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 3000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("S5").Value = 1500

# Row 6 <- old row 4 values
$ws.Range("D6").Value = 44559
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6500
$ws.Range("S6").Value = 3250

# Row 7 <- old row 5 values
$ws.Range("D7").Value = 44559
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 5000
$ws.Range("P7").Value = 5000
$ws.Range("S7").Value = 2500
